$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value, whether it is numeric-looking
# (numeric-looking text needs a leading apostrophe so Excel keeps it as text,
# then the style is reset to Normal so no stray "Text" number format sticks).
$updates = @(
    @("D2", "331.22", 1),
    @("E2", "1.06%", 1),
    @("G2", "14", 1),
    @("D3", "39.22", 1),
    @("E3", "-1.84%", 1),
    @("G3", "14", 1),
    @("D4", "5.706", 1),
    @("E4", "2.11%", 1),
    @("G4", "14", 1),
    @("D5", "0.08024", 1),
    @("E5", "-1.32%", 1),
    @("G5", "14", 1),
    @("D6", "4.494", 1),
    @("E6", "-1.21%", 1),
    @("G6", "14", 1),
    @("D7", "8.614", 1),
    @("E7", "-0.68%", 1),
    @("G7", "14", 1),
    @("D8", "1.931", 1),
    @("E8", "-2.25%", 1),
    @("G8", "14", 1),
    @("D9", "2.943", 1),
    @("E9", "-1.16%", 1),
    @("G9", "14", 1),
    @("D10", "0.9206", 1),
    @("E10", "-2.96%", 1),
    @("G10", "14", 1),
    @("D11", "0.1247", 1),
    @("E11", "-2.33%", 1),
    @("G11", "14", 1),
    @("D12", "0.1936", 1),
    @("E12", "-2.33%", 1),
    @("G12", "14", 1),
    @("D13", "8.702", 1),
    @("E13", "17.08%", 1),
    @("G13", "14", 1),
    @("D14", "0.09253", 1),
    @("E14", "0.75%", 1),
    @("G14", "14", 1),
    @("D15", "0.03565", 1),
    @("E15", "-0.19%", 1),
    @("G15", "14", 1),
    @("D16", "0.1051", 1),
    @("E16", "9.62%", 1),
    @("G16", "14", 1),
    @("E17", "-3.51%", 1),
    @("G17", "14", 1),
    @("D18", "0.006249", 1),
    @("E18", "2.75%", 1),
    @("G18", "14", 1),
    @("D19", "3.369", 1),
    @("E19", "-0.06%", 1),
    @("G19", "14", 1),
    @("D20", "0.3455", 1),
    @("E20", "-1.34%", 1),
    @("G20", "14", 1),
    @("D21", "0.1371", 1),
    @("E21", "-0.67%", 1),
    @("G21", "14", 1),
    @("D22", "0.2700", 1),
    @("E22", "8.53%", 1),
    @("G22", "14", 1),
    @("D23", "0.04440", 1),
    @("E23", "0.27%", 1),
    @("G23", "14", 1),
    @("D24", "0.001257", 1),
    @("E24", "2.64%", 1),
    @("G24", "14", 1),
    @("D25", "0.004449", 1),
    @("E25", "3.75%", 1),
    @("G25", "14", 1),
    @("D26", "0.0001204", 1),
    @("E26", "1.22%", 1),
    @("G26", "14", 1),
    @("G27", "14", 1),
    @("G28", "14", 1),
    @("G29", "14", 1),
    @("G30", "14", 1),
    @("G31", "14", 1),
    @("G32", "14", 1),
    @("G33", "14", 1),
    @("G34", "14", 1),
    @("G35", "14", 1),
    @("G36", "14", 1),
    @("G37", "14", 1),
    @("G38", "14", 1),
    @("D39", "0.02538", 1),
    @("E39", "0.97%", 1),
    @("G39", "14", 1),
    @("D40", "0.05457", 1),
    @("E40", "4.58%", 1),
    @("G40", "14", 1),
    @("D41", "0.007528", 1),
    @("E41", "-2.61%", 1),
    @("G41", "14", 1),
    @("D42", "0.009911", 1),
    @("E42", "11.63%", 1),
    @("G42", "14", 1),
    @("D43", "0.1404", 1),
    @("E43", "-1.85%", 1),
    @("G43", "14", 1),
    @("D44", "0.002115", 1),
    @("E44", "-3.42%", 1),
    @("G44", "14", 1),
    @("D45", "0.01129", 1),
    @("E45", "9.41%", 1),
    @("G45", "14", 1),
    @("D46", "0.00006821", 1),
    @("E46", "1.62%", 1),
    @("G46", "14", 1),
    @("D47", "0.00000000753", 1),
    @("E47", "0.36%", 1),
    @("G47", "14", 1),
    @("B48", "CoinbaseStockToken", 0),
    @("C48", "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin", 0),
    @("D48", "0.002287", 1),
    @("E48", "-0.56%", 1),
    @("G48", "14", 1),
    @("B49", "BOLO", 0),
    @("C49", "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo", 0),
    @("D49", "0.003070", 1),
    @("E49", "6.89%", 1),
    @("G49", "14", 1),
    @("D50", "0.00002107", 1),
    @("E50", "0.36%", 1),
    @("G50", "14", 1),
    @("D51", "0.0002007", 1),
    @("E51", "0.36%", 1),
    @("G51", "14", 1)
)

foreach ($u in $updates) {
    $addr = $u[0]
    $val = $u[1]
    $isNumeric = $u[2]
    if ($isNumeric -eq 1) {
        $ws.Range($addr).Value = "'" + $val
        $ws.Range($addr).Style = "Normal"
    } else {
        $ws.Range($addr).Value = $val
    }
}
